$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- BP-813 : Affiliate Mapping for True Independent Stations -----------
# Rename the "Affiliation Mismatch Note" header to "TrueIND"
# and the "SalesGroupName" header to "RepFirm".
$ws.Range("H1").Value = "TrueIND"
$ws.Range("J1").Value = "RepFirm"

# Bring I1/J1 formatting in line with the rest of the (bordered) header row
# by copying the format from H1, which already has the desired style.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Bring J2 and I3:J5 formatting in line with the rest of the (bordered)
# body rows by copying the format from A2.
$ws.Range("A2").Copy()
$ws.Range("J2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("A2").Copy()
$ws.Range("I3:J5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# I2 previously had no cell at all - give it the same (bordered, blank)
# formatting as its neighbours.
$ws.Range("A2").Copy()
$ws.Range("I2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
